$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 137, shifting existing rows 137..196 down to 138..197
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row 137 with the new weekly data record
$ws.Cells.Item(137, 1).Value = 9
$ws.Cells.Item(137, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(137, 3).Value = 'Metropolitana'
$ws.Cells.Item(137, 4).Value = 44529
$ws.Cells.Item(137, 5).Value = 13
$ws.Cells.Item(137, 6).Value = 'Fruta'
$ws.Cells.Item(137, 7).Value = 100103
$ws.Cells.Item(137, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(137, 9).Value = 100103001
$ws.Cells.Item(137, 10).Value = 'Cereza'
$ws.Cells.Item(137, 11).Value = 'Santina'
$ws.Cells.Item(137, 12).Value = 'Primera'
$ws.Cells.Item(137, 13).Value = 180
$ws.Cells.Item(137, 14).Value = 24000
$ws.Cells.Item(137, 15).Value = 24000
$ws.Cells.Item(137, 16).Value = 24000
$ws.Cells.Item(137, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(137, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(137, 19).Value = 2400
$ws.Cells.Item(137, 20).Value = 10
